# Break out stock.yaml completed
# 1) Fix D62:D66 on the "week" sheet so the BSE code is stored as a number
#    (it was previously written as an inline string) instead of text.
# 2) Append 5 new rows (67-71) with the newest "week" scrape for the same
#    five stocks (Astral, Glenmark, Birlasoft, Granules, Bharat Electronics).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# --- Fix up D62:D66 to be numeric values ---------------------------------
$ws.Range("D62").Value = 532830
$ws.Range("D63").Value = 532296
$ws.Range("D64").Value = 532400
$ws.Range("D65").Value = 532482
$ws.Range("D66").Value = 500049

# Helper: write a BSE code as literal text (matches the scraper's raw
# inline-string output for newly appended rows) instead of letting COM
# auto-coerce the numeric-looking string into a Number. Forcing a Text
# number format for the write then clearing formatting afterwards keeps
# the stored cell type as text while leaving the cell's style untouched
# (xfId back to the sheet default).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Append new rows 67-71 -------------------------------------------------
$ws.Range("A67").Value = 1
$ws.Range("B67").Value = "ASTRAL"
$ws.Range("C67").Value = "Astral Poly Technik Limited"
Set-TextValue $ws.Range("D67") "532830"
$ws.Range("E67").Value = 2.81
$ws.Range("F67").Value = 2273.65
$ws.Range("G67").Value = 842538
$ws.Range("H67").Value = "week"
$ws.Range("I67").Value = "21/06/2024 11:32:41"

$ws.Range("A68").Value = 2
$ws.Range("B68").Value = "GLENMARK"
$ws.Range("C68").Value = "Glenmark Pharmaceuticals Limited"
Set-TextValue $ws.Range("D68") "532296"
$ws.Range("E68").Value = -0.86
$ws.Range("F68").Value = 1230.6
$ws.Range("G68").Value = 537198
$ws.Range("H68").Value = "week"
$ws.Range("I68").Value = "21/06/2024 11:32:41"

$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "BSOFT"
$ws.Range("C69").Value = "Birlasoft Ltd"
Set-TextValue $ws.Range("D69") "532400"
$ws.Range("E69").Value = -0.96
$ws.Range("F69").Value = 681.35
$ws.Range("G69").Value = 5821357
$ws.Range("H69").Value = "week"
$ws.Range("I69").Value = "21/06/2024 11:32:41"

$ws.Range("A70").Value = 4
$ws.Range("B70").Value = "GRANULES"
$ws.Range("C70").Value = "Granules India Limited"
Set-TextValue $ws.Range("D70") "532482"
$ws.Range("E70").Value = 3.89
$ws.Range("F70").Value = 490.65
$ws.Range("G70").Value = 8523397
$ws.Range("H70").Value = "week"
$ws.Range("I70").Value = "21/06/2024 11:32:41"

$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "BEL"
$ws.Range("C71").Value = "Bharat Electronics Limited"
Set-TextValue $ws.Range("D71") "500049"
$ws.Range("E71").Value = -2.23
$ws.Range("F71").Value = 304.95
$ws.Range("G71").Value = 37620968
$ws.Range("H71").Value = "week"
$ws.Range("I71").Value = "21/06/2024 11:32:41"
